# Apply attendance count updates (0 -> 1) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where "Total Attendance Count" (D) and "Real" (E) are incremented to 1
$rowsDE = @(4, 5, 6, 9, 10, 11, 12, 13, 14, 15, 17)
foreach ($r in $rowsDE) {
    $ws.Range("D$r").Value = 1
    $ws.Range("E$r").Value = 1
}

# Rows where "Absent" (H) is incremented to 1
$rowsH = @(3, 7, 8, 16, 18)
foreach ($r in $rowsH) {
    $ws.Range("H$r").Value = 1
}

# Row 3 also has "Invalid" (G) incremented to 1
$ws.Range("G3").Value = 1

$wb.Save()
